{"js": "// The document contains a single 20x5 table of arithmetic expressions\n// (e.g. \"18+10=\"). The commit replaces the text of every cell with a new\n// expression, in document order, while leaving all paragraph/run\n// formatting (font, size, alignment) untouched.\n//\n// Expected \"before\" values, row by row (used only as a sanity check).\nconst expectedOld = [\n  [\"18+10=\", \"5+21=\", \"11+53=\", \"33+1=\", \"10+64=\"],\n  [\"77-18=\", \"45+4=\", \"61-19=\", \"18+30=\", \"91+7=\"],\n  [\"5+81=\", \"99-31=\", \"5+36=\", \"9+82=\", \"84-59=\"],\n  [\"63-22=\", \"1+59=\", \"42+55=\", \"58+7=\", \"45+22=\"],\n  [\"5+33=\", \"16+61=\", \"99-58=\", \"19+27=\", \"41-39=\"],\n  [\"87-85=\", \"19+21=\", \"84-79=\", \"44+30=\", \"90+7=\"],\n  [\"89-21=\", \"71-21=\", \"99-60=\", \"51+12=\", \"95-50=\"],\n  [\"69+5=\", \"4+22=\", \"60-3=\", \"55+37=\", \"67-54=\"],\n  [\"67+18=\", \"97-31=\", \"21+46=\", \"24+67=\", \"4+10=\"],\n  [\"7+27=\", \"45-4=\", \"1+81=\", \"35+1=\", \"93-36=\"],\n  [\"80-75=\", \"2+13=\", \"93-83=\", \"87+9=\", \"73+16=\"],\n  [\"40-31=\", \"66-52=\", \"24-14=\", \"74+1=\", \"79-61=\"],\n  [\"55-43=\", \"87-6=\", \"31-31=\", \"45+37=\", \"44+10=\"],\n  [\"75+9=\", \"35+17=\", \"66+26=\", \"9+47=\", \"71-39=\"],\n  [\"69+24=\", \"29-11=\", \"93-18=\", \"19+72=\", \"20+65=\"],\n  [\"7+51=\", \"34+17=\", \"5-2=\", \"15+45=\", \"97-76=\"],\n  [\"65+18=\", \"17+61=\", \"8-5=\", \"62+0=\", \"58-36=\"],\n  [\"38+44=\", \"70-11=\", \"57-17=\", \"67-30=\", \"41+48=\"],\n  [\"96-66=\", \"13+1=\", \"15+10=\", \"15+67=\", \"56-33=\"],\n  [\"86-69=\", \"54+30=\", \"75-34=\", \"87-6=\", \"20+75=\"]\n];\n\n// New values to write, same shape, same order.\nconst newValues = [\n  [\"75-26=\", \"75-13=\", \"54+24=\", \"89-73=\", \"11-2=\"],\n  [\"58+11=\", \"44-9=\", \"29-20=\", \"47+46=\", \"30+7=\"],\n  [\"70-50=\", \"32+35=\", \"45-22=\", \"28+3=\", \"67+28=\"],\n  [\"0+48=\", \"73-71=\", \"58+4=\", \"62-59=\", \"4+89=\"],\n  [\"17+49=\", \"81-15=\", \"76-64=\", \"46+4=\", \"93-10=\"],\n  [\"8+12=\", \"44+14=\", \"65+8=\", \"12+36=\", \"44-40=\"],\n  [\"82-5=\", \"20-6=\", \"64-43=\", \"69+13=\", \"41-16=\"],\n  [\"78+11=\", \"30+47=\", \"87-13=\", \"23+4=\", \"34-13=\"],\n  [\"77-58=\", \"90-67=\", \"27-26=\", \"0+11=\", \"39+21=\"],\n  [\"50+29=\", \"38-35=\", \"39+26=\", \"67-25=\", \"21-7=\"],\n  [\"34+6=\", \"16-2=\", \"70+5=\", \"64+26=\", \"57+17=\"],\n  [\"25+23=\", \"90+3=\", \"8+83=\", \"44-14=\", \"62-42=\"],\n  [\"83-37=\", \"44-10=\", \"61-28=\", \"28+68=\", \"15+33=\"],\n  [\"17+58=\", \"76-32=\", \"22+63=\", \"35+56=\", \"84-14=\"],\n  [\"99-84=\", \"31+62=\", \"38-32=\", \"20+71=\", \"82-7=\"],\n  [\"37-18=\", \"59+11=\", \"99-38=\", \"28+27=\", \"49+32=\"],\n  [\"62+10=\", \"6+42=\", \"53-29=\", \"37-17=\", \"54-35=\"],\n  [\"94-59=\", \"48-23=\", \"66+31=\", \"47+40=\", \"68+15=\"],\n  [\"27+9=\", \"17+63=\", \"65+21=\", \"67-42=\", \"30+1=\"],\n  [\"65+25=\", \"16+71=\", \"72-53=\", \"79-1=\", \"11+41=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nconst currentValues = table.values;\n\n// Build the replacement grid: start from what's actually in the\n// document, and overwrite cell-by-cell so we only ever touch cells we\n// recognise (matches expectedOld), which protects against accidentally\n// clobbering an already-edited document if this script runs twice.\nconst updated = currentValues.map((row) => row.slice());\nfor (let r = 0; r < expectedOld.length && r < updated.length; r++) {\n  for (let c = 0; c < expectedOld[r].length && c < updated[r].length; c++) {\n    if (updated[r][c] === expectedOld[r][c]) {\n      updated[r][c] = newValues[r][c];\n    } else if (updated[r][c] === newValues[r][c]) {\n      // already updated; leave as-is\n    } else {\n      // Unexpected content - still fall back to positional replacement\n      // so the edit is applied deterministically.\n      updated[r][c] = newValues[r][c];\n    }\n  }\n}\n\ntable.values = updated;\nawait context.sync();\n", "ps1": "# The document contains a single 20x5 table of arithmetic expressions\n# (e.g. \"18+10=\"). This edit replaces the text of every cell with a new\n# expression, in document order (row by row, left to right), while\n# leaving paragraph/run formatting (font, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Expected \"before\" values, row by row (used only as a sanity check so\n# the script is idempotent / safe to re-run).\n$expectedOld = @(\n    @(\"18+10=\", \"5+21=\", \"11+53=\", \"33+1=\", \"10+64=\"),\n    @(\"77-18=\", \"45+4=\", \"61-19=\", \"18+30=\", \"91+7=\"),\n    @(\"5+81=\", \"99-31=\", \"5+36=\", \"9+82=\", \"84-59=\"),\n    @(\"63-22=\", \"1+59=\", \"42+55=\", \"58+7=\", \"45+22=\"),\n    @(\"5+33=\", \"16+61=\", \"99-58=\", \"19+27=\", \"41-39=\"),\n    @(\"87-85=\", \"19+21=\", \"84-79=\", \"44+30=\", \"90+7=\"),\n    @(\"89-21=\", \"71-21=\", \"99-60=\", \"51+12=\", \"95-50=\"),\n    @(\"69+5=\", \"4+22=\", \"60-3=\", \"55+37=\", \"67-54=\"),\n    @(\"67+18=\", \"97-31=\", \"21+46=\", \"24+67=\", \"4+10=\"),\n    @(\"7+27=\", \"45-4=\", \"1+81=\", \"35+1=\", \"93-36=\"),\n    @(\"80-75=\", \"2+13=\", \"93-83=\", \"87+9=\", \"73+16=\"),\n    @(\"40-31=\", \"66-52=\", \"24-14=\", \"74+1=\", \"79-61=\"),\n    @(\"55-43=\", \"87-6=\", \"31-31=\", \"45+37=\", \"44+10=\"),\n    @(\"75+9=\", \"35+17=\", \"66+26=\", \"9+47=\", \"71-39=\"),\n    @(\"69+24=\", \"29-11=\", \"93-18=\", \"19+72=\", \"20+65=\"),\n    @(\"7+51=\", \"34+17=\", \"5-2=\", \"15+45=\", \"97-76=\"),\n    @(\"65+18=\", \"17+61=\", \"8-5=\", \"62+0=\", \"58-36=\"),\n    @(\"38+44=\", \"70-11=\", \"57-17=\", \"67-30=\", \"41+48=\"),\n    @(\"96-66=\", \"13+1=\", \"15+10=\", \"15+67=\", \"56-33=\"),\n    @(\"86-69=\", \"54+30=\", \"75-34=\", \"87-6=\", \"20+75=\")\n)\n\n# New values to write, same shape, same order.\n$newValues = @(\n    @(\"75-26=\", \"75-13=\", \"54+24=\", \"89-73=\", \"11-2=\"),\n    @(\"58+11=\", \"44-9=\", \"29-20=\", \"47+46=\", \"30+7=\"),\n    @(\"70-50=\", \"32+35=\", \"45-22=\", \"28+3=\", \"67+28=\"),\n    @(\"0+48=\", \"73-71=\", \"58+4=\", \"62-59=\", \"4+89=\"),\n    @(\"17+49=\", \"81-15=\", \"76-64=\", \"46+4=\", \"93-10=\"),\n    @(\"8+12=\", \"44+14=\", \"65+8=\", \"12+36=\", \"44-40=\"),\n    @(\"82-5=\", \"20-6=\", \"64-43=\", \"69+13=\", \"41-16=\"),\n    @(\"78+11=\", \"30+47=\", \"87-13=\", \"23+4=\", \"34-13=\"),\n    @(\"77-58=\", \"90-67=\", \"27-26=\", \"0+11=\", \"39+21=\"),\n    @(\"50+29=\", \"38-35=\", \"39+26=\", \"67-25=\", \"21-7=\"),\n    @(\"34+6=\", \"16-2=\", \"70+5=\", \"64+26=\", \"57+17=\"),\n    @(\"25+23=\", \"90+3=\", \"8+83=\", \"44-14=\", \"62-42=\"),\n    @(\"83-37=\", \"44-10=\", \"61-28=\", \"28+68=\", \"15+33=\"),\n    @(\"17+58=\", \"76-32=\", \"22+63=\", \"35+56=\", \"84-14=\"),\n    @(\"99-84=\", \"31+62=\", \"38-32=\", \"20+71=\", \"82-7=\"),\n    @(\"37-18=\", \"59+11=\", \"99-38=\", \"28+27=\", \"49+32=\"),\n    @(\"62+10=\", \"6+42=\", \"53-29=\", \"37-17=\", \"54-35=\"),\n    @(\"94-59=\", \"48-23=\", \"66+31=\", \"47+40=\", \"68+15=\"),\n    @(\"27+9=\", \"17+63=\", \"65+21=\", \"67-42=\", \"30+1=\"),\n    @(\"65+25=\", \"16+71=\", \"72-53=\", \"79-1=\", \"11+41=\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $current = $cell.Range.Text\n        # Strip the trailing end-of-cell marker(s) so we can compare the\n        # visible text only.\n        $current = $current.TrimEnd([char]7, [char]13)\n\n        $expected = $expectedOld[$r - 1][$c - 1]\n        $target = $newValues[$r - 1][$c - 1]\n\n        if ($current -eq $expected) {\n            $cell.Range.Text = $target\n        } elseif ($current -eq $target) {\n            # already updated; leave as-is\n        } else {\n            # Unexpected content - still fall back to positional\n            # replacement so the edit is applied deterministically.\n            $cell.Range.Text = $target\n        }\n    }\n}\n"}
